# Insert a new data row at row 182 (pushing the existing rows 182-251 down
# to 183-252, which also grows the sheet's used range from R251 to R252).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(182).Insert()

# The newly inserted row 182 is a fresh record that reuses the price/market
# data that used to live in row 182 (now shifted to row 183), but carries a
# new date and a new volume figure.
$ws.Range("A182").Value = 4
$ws.Range("B182").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C182").Value = "Los Lagos"
$ws.Range("D182").Value = 44795
$ws.Range("E182").Value = 10
$ws.Range("F182").Value = 100112039
$ws.Range("G182").Value = "Ciboulette"
$ws.Range("H182").Value = "Sin especificar"
$ws.Range("I182").Value = "Primera"
$ws.Range("J182").Value = 80
$ws.Range("K182").Value = 4000
$ws.Range("L182").Value = 4000
$ws.Range("M182").Value = 4000
$ws.Range("N182").Value = "$/docena de atados"
$ws.Range("O182").Value = "Región Metropolitana"
$ws.Range("P182").Value = 1333
$ws.Range("Q182").Value = 3
$ws.Range("R182").Value = "Hortaliza"
